$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.999.34"
$ws.Range("E2").Value = "  -0.34%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.419.59"
$ws.Range("E3").Value = "  -0.59%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.30%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "410.99"
$ws.Range("E5").Value = "  +0.02%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.48"
$ws.Range("E6").Value = "  +0.36%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.631"
$ws.Range("E7").Value = "  +0.62%  "

# Row 8
$ws.Range("E8").Value = "  +0.06%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.732"
$ws.Range("E9").Value = "  -2.21%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.139"
$ws.Range("E10").Value = "  -1.89%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "43.30"
$ws.Range("E11").Value = "  +0.42%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.20"
$ws.Range("E12").Value = "  +2.55%  "

# Row 13
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000215"
$ws.Range("E13").Value = "  +3.75%  "

# Row 14
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.957.82"
$ws.Range("E14").Value = "  -0.65%  "

# Row 15
$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.141"
$ws.Range("E15").Value = "  +0.23%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.03"
$ws.Range("E16").Value = "  -0.96%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.418.64"
$ws.Range("E17").Value = "  -1.82%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.10"
$ws.Range("E18").Value = "  +1.50%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.52"
$ws.Range("E19").Value = "  -1.14%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "61.798.27"
$ws.Range("E20").Value = "  -0.59%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "490.38"
$ws.Range("E21").Value = "  +19.77%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "92.18"
$ws.Range("E22").Value = "  +1.75%  "

# Row 23
$ws.Range("E23").Value = "  +2.75%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.49"
$ws.Range("E24").Value = "  -0.62%  "

# Row 25
$ws.Range("E25").Value = "  +3.21%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "34.25"
$ws.Range("E26").Value = "  +3.00%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.18"
$ws.Range("E27").Value = "  +6.41%  "

# Row 28
$ws.Range("E28").Value = "  -0.11%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.76"
$ws.Range("E29").Value = "  +1.48%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.70"
$ws.Range("E30").Value = "  -2.47%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.02"
$ws.Range("E31").Value = "  +0.63%  "

# Row 32
$ws.Range("E32").Value = "  -4.16%  "

# Row 33
$ws.Range("E33").Value = "  -2.45%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "41.94"
$ws.Range("E34").Value = "  -4.83%  "

# Row 35
$ws.Range("E35").Value = "  -0.08%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.64"
$ws.Range("E36").Value = "  +10.71%  "

# Row 37
$ws.Range("E37").Value = "  -2.51%  "

# Row 38
$ws.Range("E38").Value = "  +0.07%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "150.35"
$ws.Range("E39").Value = "  +6.44%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.42"
$ws.Range("E40").Value = "  +0.33%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.137"
$ws.Range("E41").Value = "  +3.62%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.323"
$ws.Range("E42").Value = "  +2.23%  "

# Row 43
$ws.Range("E43").Value = "  +0.73%  "

# Row 44
$ws.Range("E44").Value = "  +6.57%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.62"
$ws.Range("E45").Value = "  +10.20%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.26"
$ws.Range("E46").Value = "  +4.55%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.59"
$ws.Range("E47").Value = "  -1.68%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.33"
$ws.Range("E48").Value = "  +20.12%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.67"
$ws.Range("E49").Value = "  +3.22%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "117.20"
$ws.Range("E50").Value = "  +21.53%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.147"
$ws.Range("E51").Value = "  +14.64%  "
